$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The refreshed product export only carries the "Código Produto" column now,
# so wipe the old row 2-3 data across all 7 columns (B:G no longer populated).
$ws.Range("A2:G3").ClearContents()

# Write the new list of product codes coming back from the (improved) request,
# one per row, column A only.
$codes = "C-2184", "AC 30937", "HG 30784", "BD3442", "BD4190"
for ($i = 0; $i -lt $codes.Length; $i++) {
    $ws.Cells.Item(2 + $i, 1).Value = $codes[$i]
}
